$wb = $excel.ActiveWorkbook

$rushing = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

$rushing.Range("C2").Value = 13
$rushing.Range("E2").Value = 8
$rushing.Range("F2").Value = 11
$rushing.Range("C5").Value = 47
$rushing.Range("D5").Value = 24
$rushing.Range("E5").Value = 8
$rushing.Range("F5").Value = 19
$rushing.Range("C7").Value = 16
$rushing.Range("D7").Value = 10
$rushing.Range("E7").Value = 3
$rushing.Range("E8").Value = 5
$rushing.Range("C9").Value = 7
$rushing.Range("C11").Value = 2
$rushing.Range("F11").Value = 2
$rushing.Range("E12").Value = 4
$receiving.Range("C3").Value = 30
$receiving.Range("D3").Value = 22
$receiving.Range("C5").Value = 4
$receiving.Range("D5").Value = 4
$receiving.Range("C7").Value = 119
$receiving.Range("D7").Value = 94
$receiving.Range("E7").Value = 34
$receiving.Range("G7").Value = 23
$receiving.Range("E8").Value = 14
$receiving.Range("F8").Value = 6
$receiving.Range("C9").Value = 31
$receiving.Range("D9").Value = 21
$receiving.Range("E9").Value = 14
$receiving.Range("F9").Value = 10
$receiving.Range("C10").Value = 23
$receiving.Range("D10").Value = 16
$receiving.Range("E10").Value = 8
$receiving.Range("F10").Value = 3
$receiving.Range("G10").Value = 4
$receiving.Range("H10").Value = 3
$receiving.Range("E11").Value = 2
$receiving.Range("C12").Value = 7
$receiving.Range("D12").Value = 4
$receiving.Range("C13").Value = 105
$receiving.Range("D13").Value = 74
$receiving.Range("G13").Value = 16
$receiving.Range("H13").Value = 12
$receiving.Range("C14").Value = 10
$receiving.Range("D14").Value = 6
$receiving.Range("E14").Value = 2
$receiving.Range("F14").Value = 2
$receiving.Range("G14").Value = 1
$receiving.Range("H14").Value = 1
$receiving.Range("C15").Value = 7
$receiving.Range("D15").Value = 4

$rushing.Activate()
